# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet, shifting the old N/O/P columns ("Late", heading, "Outstanding")
# one place to the right, then make that sheet the active/selected sheet
# (instead of "Transactions"), matching where the author had last clicked.

$wb = $excel.ActiveWorkbook

$repayment = $wb.Worksheets.Item("Repayment Schedule")
$repayment.Columns("N").Insert()

$repayment.Activate()
$repayment.Range("P15").Select()
